
# "Generate Report for handback" - populate the per-locale handback columns
# (Latest Target File / Latest Handback File / Latest Handback DateTime) on
# the zh-cn and de-de worksheets, and flip the status text now that the
# files have been handed back.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-us"

$mdCommit = "e6a422c26b5cc56d348e9dafb91496d84a957c14"

function Set-StatusCell($ws, $addr) {
    if ($ws.Range($addr).Value2 -eq $oldStatus) {
        $ws.Range($addr).Value2 = $newStatus
    }
}

# ----- Overview sheet: flip the per-language status cells -----
$overview = $wb.Worksheets.Item("Overview")
Set-StatusCell $overview "B2"
Set-StatusCell $overview "C2"
Set-StatusCell $overview "B3"
Set-StatusCell $overview "C3"

function Apply-HandbackLocale($localeName, $xlfCommit, $handbackTime) {
    $ws = $wb.Worksheets.Item($localeName)

    # Status column flips for both tracked files.
    Set-StatusCell $ws "B2"
    Set-StatusCell $ws "B3"

    # Row 2: 0b47a6e4-381d-48f7-b1e7-fb43b3b33843.md
    $mdName2 = "0b47a6e4-381d-48f7-b1e7-fb43b3b33843.md"
    $xlfName2 = "0b47a6e4-381d-48f7-b1e7-fb43b3b33843.0ad4f96382febe82c39b58968a1c33c3b9567c2c." + $localeName + ".xlf"
    $ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/" + $mdCommit + "/e2e/" + $mdName2, "", "", $mdName2)
    $ws.Range("E2").Font.Underline = 2
    $ws.Range("E2").Font.Color = 15570276
    $ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/" + $xlfCommit + "/ol-handback/OpenLocalizationTest/oltest." + $localeName + "/yuwzho/" + $xlfName2, "", "", $xlfName2)
    $ws.Range("F2").Font.Underline = 2
    $ws.Range("F2").Font.Color = 15570276
    $ws.Range("G2").Value2 = $handbackTime

    # Row 3: 777bd7ed-b1ee-4be1-8889-e6138a54f716.md
    $mdName3 = "777bd7ed-b1ee-4be1-8889-e6138a54f716.md"
    $xlfName3 = "777bd7ed-b1ee-4be1-8889-e6138a54f716.fdbf93d3069f41d30d3f90e936dcfaa2a1779c62." + $localeName + ".xlf"
    $ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/" + $mdCommit + "/e2e/" + $mdName3, "", "", $mdName3)
    $ws.Range("E3").Font.Underline = 2
    $ws.Range("E3").Font.Color = 15570276
    $ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/" + $xlfCommit + "/ol-handback/OpenLocalizationTest/oltest." + $localeName + "/yuwzho/" + $xlfName3, "", "", $xlfName3)
    $ws.Range("F3").Font.Underline = 2
    $ws.Range("F3").Font.Color = 15570276
    $ws.Range("G3").Value2 = $handbackTime
}

Apply-HandbackLocale "zh-cn" "fdbf93d3069f41d30d3f90e936dcfaa2a1779c62" "2016-01-08 17:50:30"
Apply-HandbackLocale "de-de" "fdbf93d3069f41d30d3f90e936dcfaa2a1779c62" "2016-01-08 17:50:47"
